# Weekly price-report update: a new weekly record for "Ajo" (Chino / Primera)
# at Terminal Hortofrutícola Agro Chillán is inserted as row 122, pushing the
# existing historical rows 122-206 down to 123-207 (dimension grows to R207).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 122; this shifts rows 122..206 down to 123..207 and
# carries formatting (e.g. the date style in column D) down with them.
$ws.Rows.Item(122).Insert()

# Populate the newly inserted row 122 with the new weekly data point.
$ws.Cells.Item(122, 1).Value = 7
$ws.Cells.Item(122, 2).Value = "Terminal Hortofrutícola Agro Chillán"
$ws.Cells.Item(122, 3).Value = "Ñuble"
$ws.Cells.Item(122, 4).Value = 44634
$ws.Cells.Item(122, 5).Value = 16
$ws.Cells.Item(122, 6).Value = 100112003
$ws.Cells.Item(122, 7).Value = "Ajo"
$ws.Cells.Item(122, 8).Value = "Chino"
$ws.Cells.Item(122, 9).Value = "Primera"
$ws.Cells.Item(122, 10).Value = 60
$ws.Cells.Item(122, 11).Value = 19000
$ws.Cells.Item(122, 12).Value = 20000
$ws.Cells.Item(122, 13).Value = 19500
$ws.Cells.Item(122, 14).Value = "$/caja 10 kilos"
$ws.Cells.Item(122, 15).Value = "China"
$ws.Cells.Item(122, 16).Value = 1950
$ws.Cells.Item(122, 17).Value = 10
$ws.Cells.Item(122, 18).Value = "Hortaliza"
